$wb = $excel.ActiveWorkbook

# Physical sheet #1 (xl/worksheets/sheet1.xml) currently holds the "hotel_info"
# data; it needs to become the "review_info" sheet (header row only, A1:Y1).
$s1 = $wb.Worksheets.Item(1)
# Physical sheet #2 (xl/worksheets/sheet2.xml) currently holds the "review_info"
# header row; it needs to become the "hotel_info" sheet (headers + one data row,
# with a new "State" column inserted before "City").
$s2 = $wb.Worksheets.Item(2)

$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date",
    "review_title","review_content","review_rating","trip_month","trip_purpose",
    "value","rooms","Location","Cleanliness","Sleep Quality","Service",
    "Picture(yes=1)","respondent","response_date","response_text"
)

$hotelHeaders = @(
    "STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name",
    "English_Reviews_num","Local_Rank","Total_Reviews_num"
)

# --- Rebuild physical sheet1.xml as the new "review_info" content ---
$s1.Cells.Clear()
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $s1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Rebuild physical sheet2.xml as the new "hotel_info" content ---
$s2.Cells.Clear()
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $s2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$s2.Range("A2").Value = 43585
$s2.Range("B2").Value = "The Hotel Modern"
$s2.Range("C2").Value = "Louisiana"
$s2.Range("D2").Value = "New Orleans"
$s2.Range("E2").Value = 70130
$s2.Range("F2").Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d223122-Reviews-The_Hotel_Modern-New_Orleans_Louisiana.html"
$s2.Range("G2").Value = "The Hotel Modern"

# These three look numeric but are stored as text in the source data, so force
# text formatting before writing them, then drop back to the default style so
# no stray cell-level style survives on them.
$textRange = $s2.Range("H2:J2")
$textRange.NumberFormat = "@"
$s2.Range("H2").Value = "1002"
$s2.Range("I2").Value = "113"
$s2.Range("J2").Value = "1045"
$textRange.Style = "Normal"

# --- Fix up sheet names / tab order: name <-> sheetId now matches the target ---
$s1.Name = "temp_swap_name_zzz"
$s2.Name = "hotel_info"
$s1.Name = "review_info"
